# Add serialize/deserialize functions of CInt/CFloat.
# On the SKILL sheet, two new field columns (J = CFloat, K = CInt) are
# added to the right of the existing type-definition table (A..H),
# following the same 4-row layout: display name row (1, left blank here,
# like the other "class" style columns), field-name row (2), type row (3)
# and annotation row (4).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SKILL")

# Row 3 (type) is written before row 2 (field name) so that the shared
# string table records "CFloat" ahead of "cFloat" (and likewise "cInt"
# ahead of "CInt"), matching the order new strings were introduced.
$ws.Range("J3").Value = "CFloat"
$ws.Range("J2").Value = "cFloat"
$ws.Range("K2").Value = "cInt"
$ws.Range("K3").Value = "CInt"
$ws.Range("J4").Value = "class"
$ws.Range("K4").Value = "class"

$ws.Range("H2:H4").Copy()
$ws.Range("J2:K4").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$excel.CutCopyMode = 0

$ws.Range("K4").Select()
